$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7-12 shift down to 8-13.
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the "Trollkin" creature stat block.
$ws.Range("A7").Value = "Trollkin"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = "Darkness"
$ws.Range("Q7").Value = 40
$ws.Range("T7").Value = 6
$ws.Range("Z7").Value = "Humanoid"

# Match the author's final UI selection state recorded in the saved file.
[void]$ws.Range("Q25").Select()
